$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.690.80'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').Value = '3.146.23'
$ws.Range('E3').Value = '  +3.32%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.18'
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.73'
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('E7').Value = '  +30.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.369'
$ws.Range('E8').Value = '  +1.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = '3.145.10'
$ws.Range('E10').Value = '  +3.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.758'
$ws.Range('E11').Value = '  +14.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.203'
$ws.Range('E12').Value = '  +8.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.67'
$ws.Range('E13').Value = '  +5.79%  '
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.87'
$ws.Range('E15').Value = '  +8.93%  '
$ws.Range('D16').Value = '90.339.89'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('D17').Value = '3.724.19'
$ws.Range('E17').Value = '  +3.72%  '
$ws.Range('D18').Value = '3.127.98'
$ws.Range('E18').Value = '  +3.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.66'
$ws.Range('E19').Value = '  +9.35%  '
$ws.Range('E20').Value = '  +6.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '463.15'
$ws.Range('E21').Value = '  +8.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0000209'
$ws.Range('E22').Value = '  -2.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.07'
$ws.Range('E23').Value = '  +10.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.27'
$ws.Range('E24').Value = '  +5.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.71'
$ws.Range('E25').Value = '  +5.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '91.16'
$ws.Range('E26').Value = '  +9.05%  '
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('D28').Value = '3.314.55'
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.10'
$ws.Range('E32').Value = '  +11.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.19'
$ws.Range('E33').Value = '  +18.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '516.06'
$ws.Range('E34').Value = '  +2.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.182'
$ws.Range('E35').Value = '  +32.66%  '
$ws.Range('E36').Value = '  +7.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.60'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.85'
$ws.Range('E38').Value = '  +3.18%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.141'
$ws.Range('E39').Value = '  +7.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.29'
$ws.Range('E40').Value = '  +4.48%  '
$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.20'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0853'
$ws.Range('E42').Value = '  +25.40%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.412'
$ws.Range('E44').Value = '  +13.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.96'
$ws.Range('E45').Value = '  +7.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '149.94'
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.64'
$ws.Range('E48').Value = '  +5.40%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.55'
$ws.Range('E49').Value = '  +11.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.35'
$ws.Range('E50').Value = '  +10.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.668'
$ws.Range('E51').Value = '  +13.87%  '
